# Update the email addresses in the data list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "nityaranjn7897756843@gmail.com"
$ws.Range("C3").Value = "abh090653646@gmail.com"

# Leave the cursor where it was when the workbook was last saved.
$ws.Range("C13").Select()
